# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# This overwrites the literal values in column G (rows 2-35) of the active
# worksheet with the newly-computed "K" values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(4, 3, 6, 5, 9, 7, 5, 4, 7, 10, 7, 9, 4, 8, 9, 4, 6, 8, 9, 5, 7, 4, 2, 9, 7, 2, 10, 11, 6, 4, 5, 4, 5, 3)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
